$d = $word.ActiveDocument

# --- 1. Qualifications bullet list: set font size to 10pt (sz=20) ---
# Paragraphs 4-8 are the five "Qualifications" bullets. Setting Font.Size
# on the paragraph Range applies <w:sz w:val="20"/> to both the paragraph
# mark rPr and every run's rPr, matching the diff.
for ($i = 4; $i -le 8; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Size = 10
}

# --- 2. Objective paragraph: only the sentence after "Objective<tab>" gets sz=20 ---
$objP = $d.Paragraphs.Item(10)
$objText = $objP.Range.Text
$tabIdx = $objText.IndexOf([char]9)
$bodyStart = $objP.Range.Start + $tabIdx + 1
$bodyRange = $d.Range($bodyStart, $objP.Range.End)
$bodyRange.Font.Size = 10

# --- 3. Move the "_GoBack" bookmark from the end of the document to inside
#        the "Designed a matching card game..." run (splits that run in two). ---
$d.Bookmarks.Item("_GoBack").Delete()

$cardP = $d.Paragraphs.Item(24)
$splitPos = $cardP.Range.Start + 10
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)
